$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new parameter rows above the old row 14 (pushes everything that
# used to be rows 14-34 down to rows 16-36).
# ---------------------------------------------------------------------------
$ws.Range("A14:A15").EntireRow.Insert()

# Pick up the same look as the rest of the "plain" parameter rows (e.g. row 13)
# for the two freshly inserted rows.
$ws.Range("A13:B13").Copy()
$ws.Range("A14:B15").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D13").Copy()
$ws.Range("D14:D15").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New row 14: init_pr_ever_talking_therapy_if_diagnosed
$ws.Range("A14").Value2 = "init_pr_ever_talking_therapy_if_diagnosed"
$ws.Range("B14").Value2 = 1
$ws.Range("D14").Value2 = "We assume that talking therapy happens as part of diagnosis"

# New row 15: init_pr_ever_self_harmed_if_ever_depr
$ws.Range("A15").Value2 = "init_pr_ever_self_harmed_if_ever_depr"
$ws.Range("B15").Value2 = 0.004
$ws.Range("D15").Value2 = "consistent with rate of incident self harm "

# ---------------------------------------------------------------------------
# The workbook was resaved by a newer Excel build with a slightly larger
# default row height / font metric, which re-flows the autofit row heights
# of every wrapped-text row below.  Re-apply autofit everywhere, then pin the
# handful of rows that need an explicit (wrapped) height.
# ---------------------------------------------------------------------------
$plainRows = @(17,20,22,23,24,25,26,27,28,29,30,32,33,34,35,36)
foreach ($r in $plainRows) {
    $ws.Rows.Item($r).AutoFit()
}

$ws.Rows.Item(16).RowHeight = 62
$ws.Rows.Item(18).RowHeight = 46.5
$ws.Rows.Item(19).RowHeight = 46.5
$ws.Rows.Item(21).RowHeight = 31
$ws.Rows.Item(31).RowHeight = 62

# ---------------------------------------------------------------------------
# Selection / scroll position, matching the saved view state.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("A16").Select()

Write-Host "applied depression resourcefile update"
